$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.559823036193848
$ws.Range("B1").Value = 2.838309526443481
$ws.Range("C1").Value = 3.979631900787354
$ws.Range("D1").Value = 1.517715573310852
$ws.Range("E1").Value = 0.9957099556922913
